# Update the daily 'cryptos' price/volume snapshot (refreshed by the scraper).
# Only specific cells change: Price (D) / Volume(1h) (E) figures are refreshed,
# and two coin pairs (rows 26/27 and 45/46) swap rank position, which also
# swaps their Coin name / Link / Price / Volume values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '41.832.09'
$ws.Cells.Item(2, 5).Value = '  +4.08%  '
# Row 3
$ws.Cells.Item(3, 4).Value = '2.273.72'
$ws.Cells.Item(3, 5).Value = '  +2.12%  '
# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.01%  '
# Row 5
$ws.Cells.Item(5, 4).Value = '''304.64'
$ws.Cells.Item(5, 5).Value = '  +3.60%  '
# Row 6
$ws.Cells.Item(6, 4).Value = '''93.11'
$ws.Cells.Item(6, 5).Value = '  +5.96%  '
# Row 7
$ws.Cells.Item(7, 5).Value = '  +3.78%  '
# Row 8
$ws.Cells.Item(8, 5).Value = '  -0.01%  '
# Row 9
$ws.Cells.Item(9, 5).Value = '  +3.56%  '
# Row 10
$ws.Cells.Item(10, 4).Value = '''32.71'
$ws.Cells.Item(10, 5).Value = '  +6.06%  '
# Row 11
$ws.Cells.Item(11, 4).Value = '''53.86'
$ws.Cells.Item(11, 5).Value = '  +5.63%  '
# Row 12
$ws.Cells.Item(12, 4).Value = '''0.0803'
# Row 13
$ws.Cells.Item(13, 5).Value = '  +1.98%  '
# Row 14
$ws.Cells.Item(14, 4).Value = '''6.69'
$ws.Cells.Item(14, 5).Value = '  +3.61%  '
# Row 15
$ws.Cells.Item(15, 4).Value = '2.626.66'
$ws.Cells.Item(15, 5).Value = '  +1.58%  '
# Row 16
$ws.Cells.Item(16, 4).Value = '''14.27'
$ws.Cells.Item(16, 5).Value = '  +2.82%  '
# Row 17
$ws.Cells.Item(17, 4).Value = '2.284.68'
$ws.Cells.Item(17, 5).Value = '  +2.29%  '
# Row 18
$ws.Cells.Item(18, 4).Value = '''0.763'
$ws.Cells.Item(18, 5).Value = '  +3.37%  '
# Row 19
$ws.Cells.Item(19, 4).Value = '41.785.58'
$ws.Cells.Item(19, 5).Value = '  +4.15%  '
# Row 20
$ws.Cells.Item(20, 4).Value = '''12.31'
$ws.Cells.Item(20, 5).Value = '  +9.05%  '
# Row 21
$ws.Cells.Item(21, 5).Value = '  +2.35%  '
# Row 22
$ws.Cells.Item(22, 5).Value = '  +2.84%  '
# Row 23
$ws.Cells.Item(23, 4).Value = '''67.36'
$ws.Cells.Item(23, 5).Value = '  +2.44%  '
# Row 24
$ws.Cells.Item(24, 4).Value = '''243.84'
$ws.Cells.Item(24, 5).Value = '  +2.94%  '
# Row 25
$ws.Cells.Item(25, 4).Value = '''2.60'
$ws.Cells.Item(25, 5).Value = '  +4.68%  '
# Row 26
$ws.Cells.Item(26, 2).Value = 'ImmutableX'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(26, 4).Value = '''1.93'
$ws.Cells.Item(26, 5).Value = '  +5.35%  '
# Row 27
$ws.Cells.Item(27, 2).Value = 'Dai'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(27, 4).Value = '''1.00'
$ws.Cells.Item(27, 5).Value = '  +0.04%  '
# Row 28
$ws.Cells.Item(28, 4).Value = '''24.24'
$ws.Cells.Item(28, 5).Value = '  +3.96%  '
# Row 29
$ws.Cells.Item(29, 4).Value = '''9.62'
$ws.Cells.Item(29, 5).Value = '  +3.19%  '
# Row 30
$ws.Cells.Item(30, 5).Value = '  +0.93%  '
# Row 31
$ws.Cells.Item(31, 4).Value = '''34.15'
$ws.Cells.Item(31, 5).Value = '  +6.96%  '
# Row 32
$ws.Cells.Item(32, 4).Value = '''158.64'
$ws.Cells.Item(32, 5).Value = '  +0.34%  '
# Row 33
$ws.Cells.Item(33, 5).Value = '  +0.00%  '
# Row 34
$ws.Cells.Item(34, 5).Value = '  +4.12%  '
# Row 35
$ws.Cells.Item(35, 4).Value = '''0.0752'
$ws.Cells.Item(35, 5).Value = '  +4.79%  '
# Row 36
$ws.Cells.Item(36, 4).Value = '''3.05'
$ws.Cells.Item(36, 5).Value = '  +0.13%  '
# Row 37
$ws.Cells.Item(37, 5).Value = '  +3.19%  '
# Row 38
$ws.Cells.Item(38, 5).Value = '  +7.68%  '
# Row 39
$ws.Cells.Item(39, 5).Value = '  +3.00%  '
# Row 40
$ws.Cells.Item(40, 5).Value = '  +5.03%  '
# Row 41
$ws.Cells.Item(41, 4).Value = '''1.82'
$ws.Cells.Item(41, 5).Value = '  +3.86%  '
# Row 42
$ws.Cells.Item(42, 5).Value = '  +5.47%  '
# Row 43
$ws.Cells.Item(43, 4).Value = '2.069.40'
$ws.Cells.Item(43, 5).Value = '  -0.70%  '
# Row 44
$ws.Cells.Item(44, 4).Value = '''19.90'
$ws.Cells.Item(44, 5).Value = '  +4.64%  '
# Row 45
$ws.Cells.Item(45, 2).Value = 'VeChain'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(45, 4).Value = '''0.0279'
$ws.Cells.Item(45, 5).Value = '  +3.09%  '
# Row 46
$ws.Cells.Item(46, 2).Value = 'FraxShare'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(46, 4).Value = '''10.40'
$ws.Cells.Item(46, 5).Value = '  +2.59%  '
# Row 47
$ws.Cells.Item(47, 4).Value = '''2.91'
$ws.Cells.Item(47, 5).Value = '  +6.01%  '
# Row 48
$ws.Cells.Item(48, 5).Value = '  +4.64%  '
# Row 49
$ws.Cells.Item(49, 5).Value = '  +2.83%  '
# Row 50
$ws.Cells.Item(50, 4).Value = '''72.87'
$ws.Cells.Item(50, 5).Value = '  +7.10%  '
# Row 51
$ws.Cells.Item(51, 5).Value = '  +3.65%  '
